$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final table (header + 7 data rows, now including a new "price" column F)
$data = @(
    @("product code", "name", "size", "amount", "color", "price"),
    @(1, "shirts", "s/m/l", 7, "black", 40),
    @(2, "jeans", "s/m/l", 10, "black", 50),
    @(3, "shoes", "s/m/l", 10, "black", 60),
    @(4, "coats", "s/m/l", 10, "black", 24),
    @(5, "belts", "s/m/l", 34, "red", 12),
    @(6, "shirt", "s/m/l", 6, "white", 45),
    @(7, "shorts", "s/m/l", 10, "blue", 70)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}

$ws.Range("A1:F8").Select()
